{"js": "// Remove the hyphen from \"parenting-time\" -> \"parenting time\"\n// (commit: \"Removed hyphen from word 'parenting-time' in attachment\")\nconst body = context.document.body;\n\nconst results = body.search(\"parenting-time\", { matchCase: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"parenting time\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Remove the hyphen from \"parenting-time\" -> \"parenting time\"\n# (commit: \"Removed hyphen from word 'parenting-time' in attachment\")\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"parenting-time\"\n$find.Replacement.Text = \"parenting time\"\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $false\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\n\n# wdReplace: 0 = wdReplaceNone, 1 = wdReplaceOne, 2 = wdReplaceAll\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
